# Updated CHE_grids model - 2025-08-09 01:22
# The "solar" sheet's grid_cell column (AG) values for rows 7-22 and 26
# are re-pointed to a different CHE_<n> label (a re-shuffle of the grid
# cell assignments), while every other column/row is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$updates = @{
    7  = "CHE_20"
    8  = "CHE_1"
    9  = "CHE_6"
    10 = "CHE_10"
    11 = "CHE_22"
    12 = "CHE_12"
    13 = "CHE_17"
    14 = "CHE_19"
    15 = "CHE_0"
    16 = "CHE_11"
    17 = "CHE_15"
    18 = "CHE_25"
    19 = "CHE_14"
    20 = "CHE_18"
    21 = "CHE_3"
    22 = "CHE_13"
    26 = "CHE_7"
}

foreach ($row in $updates.Keys) {
    $ws.Range("AG$row").Value2 = $updates[$row]
}

$wb.Save()
